$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.372367262840271
$ws.Range("B1").Value = 0.261387437582016
$ws.Range("C1").Value = 0.4040209054946899
$ws.Range("D1").Value = 4.334035396575928
$ws.Range("E1").Value = 2.333630323410034
